$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 18 for the "0.3333333333333333C" usable capacity ratio,
# just below the existing "0.5C" row (old row 17).
$ws.Rows.Item(18).Insert()

# Insert two new rows at 21-22 for the "0.16666666666666666C" and "0.125C"
# usable capacity ratios, just below the existing "0.2C" row (which has now
# shifted down to row 20 because of the insert above).
$ws.Range("A21:A22").EntireRow.Insert()

# Fill in the newly inserted row 18.
$ws.Range("A18").Value = "usable_capacity_ratio_at_0.3333333333333333C"
$ws.Range("B18").Value = 0.92
$ws.Range("C18").Value = 0.93
$ws.Range("D18").Value = 0.97

# Fill in the newly inserted rows 21 and 22.
$ws.Range("A21").Value = "usable_capacity_ratio_at_0.16666666666666666C"
$ws.Range("B21").Value = 0.95
$ws.Range("C21").Value = 0.96
$ws.Range("D21").Value = 0.97

$ws.Range("A22").Value = "usable_capacity_ratio_at_0.125C"
$ws.Range("B22").Value = 0.95
$ws.Range("C22").Value = 0.96
$ws.Range("D22").Value = 0.97

# Widen column A to fit the new, longer labels.
$ws.Columns.Item(1).ColumnWidth = 42.85

# Update the active selection to reflect where editing left off.
$ws.Range("D23").Select() | Out-Null
